$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.9521773333333332
$ws.Range("H2").Value = 2.856532
$ws.Range("I2").Value = 0.005888173680712573
$ws.Range("J2").Value = 0.005888173680712573
$ws.Range("M2").Value = 188.5745136666667
$ws.Range("N2").Value = 565.723541
$ws.Range("O2").Value = 0.3378563924127341
$ws.Range("P2").Value = 0.3378563924127341
$ws.Range("Q2").Value = 179.5563775577569
$ws.Range("R2").Value = 1616.007398019812
$ws.Range("S2").Value = 0.00198935711766516
$ws.Range("T2").Value = 0.00198935711766516
$ws.Range("G3").Value = 0.9521773333333332
$ws.Range("H3").Value = 2.856532
$ws.Range("I3").Value = 0.005888173680712573
$ws.Range("J3").Value = 0.005888173680712573
$ws.Range("O3").Value = 0.3573715899422427
$ws.Range("P3").Value = 0.3573715899422427
$ws.Range("Q3").Value = 189.9278793390284
$ws.Range("R3").Value = 1709.350914051256
$ws.Range("S3").Value = 0.00210426599013232
$ws.Range("T3").Value = 0.00210426599013232
$ws.Range("G4").Value = 0.9521773333333332
$ws.Range("H4").Value = 2.856532
$ws.Range("I4").Value = 0.005888173680712573
$ws.Range("J4").Value = 0.005888173680712573
$ws.Range("M4").Value = 125.2744573333333
$ws.Range("N4").Value = 375.823372
$ws.Range("O4").Value = 0.224445898828716
$ws.Range("P4").Value = 0.224445898828716
$ws.Range("Q4").Value = 119.2834987184338
$ws.Range("R4").Value = 1073.551488465904
$ws.Range("S4").Value = 0.001321576434227122
$ws.Range("T4").Value = 0.001321576434227122
$ws.Range("G5").Value = 0.9521773333333332
$ws.Range("H5").Value = 2.856532
$ws.Range("I5").Value = 0.005888173680712573
$ws.Range("J5").Value = 0.005888173680712573
$ws.Range("M5").Value = 44.83401566666667
$ws.Range("N5").Value = 134.502047
$ws.Range("O5").Value = 0.08032611881630715
$ws.Range("P5").Value = 0.08032611881630715
$ws.Range("Q5").Value = 42.68993348011155
$ws.Range("R5").Value = 384.209401321004
$ws.Range("S5").Value = 0.0004729741386879708
$ws.Range("T5").Value = 0.0004729741386879708
$ws.Range("I6").Value = 0.5421346526196088
$ws.Range("J6").Value = 0.5421346526196088
$ws.Range("M6").Value = 188.5745136666667
$ws.Range("N6").Value = 565.723541
$ws.Range("O6").Value = 0.3378563924127341
$ws.Range("P6").Value = 0.3378563924127341
$ws.Range("Q6").Value = 16532.07592903909
$ws.Range("R6").Value = 148788.6833613518
$ws.Range("S6").Value = 0.1831636579359918
$ws.Range("T6").Value = 0.1831636579359918
$ws.Range("I7").Value = 0.5421346526196088
$ws.Range("J7").Value = 0.5421346526196088
$ws.Range("O7").Value = 0.3573715899422427
$ws.Range("P7").Value = 0.3573715899422427
$ws.Range("S7").Value = 0.193743522769455
$ws.Range("T7").Value = 0.193743522769455
$ws.Range("I8").Value = 0.5421346526196088
$ws.Range("J8").Value = 0.5421346526196088
$ws.Range("M8").Value = 125.2744573333333
$ws.Range("N8").Value = 375.823372
$ws.Range("O8").Value = 0.224445898828716
$ws.Range("P8").Value = 0.224445898828716
$ws.Range("Q8").Value = 10982.64447477094
$ws.Range("R8").Value = 98843.8002729385
$ws.Range("S8").Value = 0.1216798993934018
$ws.Range("T8").Value = 0.1216798993934018
$ws.Range("I9").Value = 0.5421346526196088
$ws.Range("J9").Value = 0.5421346526196088
$ws.Range("M9").Value = 44.83401566666667
$ws.Range("N9").Value = 134.502047
$ws.Range("O9").Value = 0.08032611881630715
$ws.Range("P9").Value = 0.08032611881630715
$ws.Range("Q9").Value = 3930.538315030423
$ws.Range("R9").Value = 35374.84483527381
$ws.Range("S9").Value = 0.0435475725207601
$ws.Range("T9").Value = 0.0435475725207601
$ws.Range("G10").Value = 72.96496833333333
$ws.Range("H10").Value = 218.894905
$ws.Range("I10").Value = 0.4512083948168896
$ws.Range("J10").Value = 0.4512083948168896
$ws.Range("M10").Value = 188.5745136666667
$ws.Range("N10").Value = 565.723541
$ws.Range("O10").Value = 0.3378563924127341
$ws.Range("P10").Value = 0.3378563924127341
$ws.Range("Q10").Value = 13759.33341816207
$ws.Range("R10").Value = 123834.0007634586
$ws.Range("S10").Value = 0.1524436404991749
$ws.Range("T10").Value = 0.1524436404991749
$ws.Range("G11").Value = 72.96496833333333
$ws.Range("H11").Value = 218.894905
$ws.Range("I11").Value = 0.4512083948168896
$ws.Range("J11").Value = 0.4512083948168896
$ws.Range("O11").Value = 0.3573715899422427
$ws.Range("P11").Value = 0.3573715899422427
$ws.Range("Q11").Value = 14554.09745270422
$ws.Range("R11").Value = 130986.877074338
$ws.Range("S11").Value = 0.161249061450999
$ws.Range("T11").Value = 0.161249061450999
$ws.Range("G12").Value = 72.96496833333333
$ws.Range("H12").Value = 218.894905
$ws.Range("I12").Value = 0.4512083948168896
$ws.Range("J12").Value = 0.4512083948168896
$ws.Range("M12").Value = 125.2744573333333
$ws.Range("N12").Value = 375.823372
$ws.Range("O12").Value = 0.224445898828716
$ws.Range("P12").Value = 0.224445898828716
$ws.Range("Q12").Value = 9140.646812302184
$ws.Range("R12").Value = 82265.82131071966
$ws.Range("S12").Value = 0.1012718737337389
$ws.Range("T12").Value = 0.1012718737337389
$ws.Range("G13").Value = 72.96496833333333
$ws.Range("H13").Value = 218.894905
$ws.Range("I13").Value = 0.4512083948168896
$ws.Range("J13").Value = 0.4512083948168896
$ws.Range("M13").Value = 44.83401566666667
$ws.Range("N13").Value = 134.502047
$ws.Range("O13").Value = 0.08032611881630715
$ws.Range("P13").Value = 0.08032611881630715
$ws.Range("Q13").Value = 3271.312533374504
$ws.Range("R13").Value = 29441.81280037054
$ws.Range("S13").Value = 0.0362438191329767
$ws.Range("T13").Value = 0.0362438191329767
$ws.Range("G14").Value = 0.1243193333333333
$ws.Range("H14").Value = 0.372958
$ws.Range("I14").Value = 0.0007687788827890604
$ws.Range("J14").Value = 0.0007687788827890604
$ws.Range("M14").Value = 188.5745136666667
$ws.Range("N14").Value = 565.723541
$ws.Range("O14").Value = 0.3378563924127341
$ws.Range("P14").Value = 0.3378563924127341
$ws.Range("Q14").Value = 23.44345782269756
$ws.Range("R14").Value = 210.991120404278
$ws.Range("S14").Value = 0.0002597368599022041
$ws.Range("T14").Value = 0.0002597368599022041
$ws.Range("G15").Value = 0.1243193333333333
$ws.Range("H15").Value = 0.372958
$ws.Range("I15").Value = 0.0007687788827890604
$ws.Range("J15").Value = 0.0007687788827890604
$ws.Range("O15").Value = 0.3573715899422427
$ws.Range("P15").Value = 0.3573715899422427
$ws.Range("Q15").Value = 24.79759443357378
$ws.Range("R15").Value = 223.178349902164
$ws.Range("S15").Value = 0.0002747397316563476
$ws.Range("T15").Value = 0.0002747397316563476
$ws.Range("G16").Value = 0.1243193333333333
$ws.Range("H16").Value = 0.372958
$ws.Range("I16").Value = 0.0007687788827890604
$ws.Range("J16").Value = 0.0007687788827890604
$ws.Range("M16").Value = 125.2744573333333
$ws.Range("N16").Value = 375.823372
$ws.Range("O16").Value = 0.224445898828716
$ws.Range("P16").Value = 0.224445898828716
$ws.Range("Q16").Value = 15.57403701937511
$ws.Range("R16").Value = 140.166333174376
$ws.Range("S16").Value = 0.0001725492673481268
$ws.Range("T16").Value = 0.0001725492673481267
$ws.Range("G17").Value = 0.1243193333333333
$ws.Range("H17").Value = 0.372958
$ws.Range("I17").Value = 0.0007687788827890604
$ws.Range("J17").Value = 0.0007687788827890604
$ws.Range("M17").Value = 44.83401566666667
$ws.Range("N17").Value = 134.502047
$ws.Range("O17").Value = 0.08032611881630715
$ws.Range("P17").Value = 0.08032611881630715
$ws.Range("Q17").Value = 5.573734938336222
$ws.Range("R17").Value = 50.16361444502601
$ws.Range("S17").Value = 0.00006175302388238193
$ws.Range("T17").Value = 0.00006175302388238193
